# Updates the "想去人数" (F column) values on the "展览" and "全部类型"
# worksheets to reflect newly generated output data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> (row -> new F value)
$updates = @{
    "展览" = @{
        3  = 577
        7  = 14480
        8  = 398
        9  = 653
        10 = 15116
        12 = 8532
        15 = 53
        16 = 172
        21 = 4
        23 = 1065
        25 = 5
        26 = 45
        30 = 17
        31 = 20
        32 = 226
        33 = 250
        34 = 408
        35 = 108
        36 = 5245
        37 = 5225
    }
    "全部类型" = @{
        3  = 577
        7  = 14480
        8  = 398
        9  = 653
        10 = 15116
        12 = 8532
        16 = 53
        17 = 172
        22 = 4
        24 = 1065
        26 = 5
        27 = 45
        33 = 17
        34 = 20
        35 = 226
        36 = 250
        37 = 408
        38 = 108
        39 = 5245
        40 = 5225
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
